# Auto-generated edit script: updates crypto price/volume data
# per the commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.467.73"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "1.848.34"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'243.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.6585"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.2992"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.07481"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'24.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07640"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.842.22"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.020"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6840"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'83.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.000009519"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'6.140"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "29.507.88"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.074.95"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'236.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'12.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'7.648"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.08%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value = "'0.1423"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'156.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'8.487"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'17.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.06026"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "'1.488"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.250"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.135"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.075"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.179"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.851"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7225"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.596"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.802"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01780"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.197.22"
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.236"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.9087"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'0.9998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.010.74"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'101.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'66.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.438"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.60%  "
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").Value = "'0.4058"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000117"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.033"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.653"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.84%  "
